$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 66, shifting existing rows 66-72 down to 67-73
$ws.Rows("66:66").Insert()

# Populate the newly inserted row 66 with the new data point
$ws.Range("A66").Value = 9
$ws.Range("B66").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C66").Value = "Metropolitana"
$ws.Range("D66").Value = 44491
$ws.Range("E66").Value = 13
$ws.Range("F66").Value = 100112005
$ws.Range("G66").Value = "Puerro"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 160
$ws.Range("K66").Value = 7000
$ws.Range("L66").Value = 8000
$ws.Range("M66").Value = 7500
$ws.Range("N66").Value = "$/paquete 20 unidades"
$ws.Range("O66").Value = "Provincia de Chacabuco"
$ws.Range("P66").Value = 375
$ws.Range("Q66").Value = 20
$ws.Range("R66").Value = "Hortaliza"

# Match the style of the new Date cell (D) to the other D cells
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat
